$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (e.g. H1) by copying its formatting over.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new data columns I (I0) and J (IF) for rows 2-10
$values = @(
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(9, 10),
    @(3, 5),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
